$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The function for no decision has been successfully called.`n"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision about Friday's movie could not be finalized, and therefore no selection was made.`n"
$ws.Range("D5").Value = "no_decision, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie for Friday, resulting in no decision.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be acquired for the Friday showing.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision-making process resulted in no agreement on which movie to show on Friday.`n"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" have been acquired for the movie to be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday was made during the discussion.`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision process did not reach a conclusion, so there is no plan for Friday's movie.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: I have recorded the decision as no_decision, indicating that the committee could not finalize a choice for Friday's movie.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made about Friday's movie.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, and there will be no movie acquired for Friday.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for showing on Friday.`n"
$ws.Range("D24").Value = "both_movies, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding which movie to show on Friday, so no selection has been made.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement on a movie for Friday.`n"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire was not made.`n"
$ws.Range("D30").Value = "no_decision, "
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("D31").Value = "both_movies, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The rights to both movies, `"Oppenheimer`" and `"Barbie,`" have been acquired.`n"
$ws.Range("D33").Value = "both_movies, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D34").Value = "both_movies, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for the movie `"Barbie`" to be shown on Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement.`n"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the screening.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision was made to not select a movie for Friday.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to `"Oppenheimer`" have been acquired.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie acquired for Friday's showing.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been made to call the no_decision function, indicating that no consensus was reached regarding the movie to be shown on Friday.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision process concluded without reaching a consensus, so no further action will be taken regarding the movie selection for Friday.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" have been acquired.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: None`n`nMSG: No decision has been made regarding which movie to show on Friday.`n"
$ws.Range("D50").Value = "no_decision, , no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights to `"Barbie`" for the movie to be shown on Friday.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not made, and thus the outcome is classified as `"no decision.`"`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been successfully recorded.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision to acquire a movie for Friday could not be reached, resulting in no selection being made.`n"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None`n`nMSG: The committee did not arrive at a decision regarding which movie to show on Friday.`n"
$ws.Range("D57").Value = "no_decision, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("D59").Value = "no_decision, "

Write-Host "Applied classification updates"